$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "306.03"
Set-TextValue 2 5 "1.56%"
Set-TextValue 3 4 "36.35"
Set-TextValue 3 5 "-0.58%"
Set-TextValue 4 4 "5.064"
Set-TextValue 4 5 "1.59%"
Set-TextValue 5 4 "0.07927"
Set-TextValue 5 5 "3.02%"
Set-TextValue 6 4 "2.181"
Set-TextValue 6 5 "4.21%"
Set-TextValue 7 4 "4.173"
Set-TextValue 7 5 "3.60%"
Set-TextValue 8 4 "8.015"
Set-TextValue 8 5 "1.24%"
Set-TextValue 9 4 "0.9299"
Set-TextValue 9 5 "1.72%"
Set-TextValue 10 4 "0.09835"
Set-TextValue 10 5 "1.98%"
Set-TextValue 11 4 "0.1876"
Set-TextValue 11 5 "0.70%"
Set-TextValue 12 4 "0.09095"
Set-TextValue 12 5 "6.70%"
Set-TextValue 13 4 "0.03711"
Set-TextValue 13 5 "4.86%"
Set-TextValue 14 4 "0.09915"
Set-TextValue 14 5 "-0.65%"
Set-TextValue 15 4 "0.001433"
Set-TextValue 15 5 "-3.26%"
Set-TextValue 16 4 "0.005634"
Set-TextValue 16 5 "-0.24%"
Set-TextValue 17 4 "3.462"
Set-TextValue 17 5 "-0.07%"
Set-TextValue 18 5 "8.11%"
Set-TextValue 19 5 "-0.37%"
Set-TextValue 20 4 "0.1336"
Set-TextValue 20 5 "0.74%"
Set-TextValue 21 4 "5.118"
Set-TextValue 21 5 "7.56%"
Set-TextValue 22 4 "0.2188"
Set-TextValue 22 5 "-0.56%"
Set-TextValue 23 4 "0.04551"
Set-TextValue 23 5 "-1.07%"
Set-TextValue 24 5 "0.59%"
Set-TextValue 25 4 "0.004783"
Set-TextValue 25 5 "-6.23%"
Set-TextValue 26 4 "0.0001301"
Set-TextValue 26 5 "-7.13%"
Set-TextValue 39 4 "0.01923"
Set-TextValue 39 5 "9.24%"
Set-TextValue 40 4 "0.04934"
Set-TextValue 40 5 "7.14%"
Set-TextValue 41 4 "0.007752"
Set-TextValue 41 5 "2.77%"
Set-TextValue 42 4 "0.1396"
Set-TextValue 42 5 "0.49%"
Set-TextValue 43 4 "0.007803"
Set-TextValue 43 5 "1.05%"
Set-TextValue 44 4 "0.002112"
Set-TextValue 44 5 "-5.75%"
Set-TextValue 45 4 "0.01124"
Set-TextValue 45 5 "9.05%"
Set-TextValue 46 4 "0.00006226"
Set-TextValue 46 5 "1.83%"
Set-TextValue 47 4 "0.00000000750"
Set-TextValue 47 5 "-0.09%"
Set-TextValue 48 4 "52.26"
Set-TextValue 48 5 "50.82%"
Set-TextValue 49 4 "0.001799"
Set-TextValue 49 5 "-10.07%"
Set-TextValue 50 4 "0.00002099"
Set-TextValue 50 5 "-0.09%"
Set-TextValue 51 4 "0.0001999"
Set-TextValue 51 5 "-0.09%"
